$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# ---------------------------------------------------------------------------
# 1. New "Login" column (G) on the user sheet.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Login"
# Give the new header cell the same yellow fill as the other header cells but
# without the left/top alignment override (matches the new style produced by
# the real edit).
$ws.Range("G1").Interior.Color = 65535

$ws.Columns.Item(7).ColumnWidth = 10.73

# ---------------------------------------------------------------------------
# 2. Promote row 7 (previously just the stray "testaccount6@mail.com" cell)
#    into a full data row, and add new users (rows 8-10).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Alf"
$ws.Range("C7").Value = "Doe"
$ws.Range("D7").Value = "testaccount6@mail.com"
$ws.Range("E7").Value = "Tester123@"
$ws.Range("F7").Value = "N"
$ws.Range("G7").Value = "N"
$ws.Range("A7").Style = $ws.Range("A6").Style
$ws.Range("B7").Style = $ws.Range("B6").Style
$ws.Range("C7").Style = $ws.Range("C6").Style
$ws.Range("E7").Style = $ws.Range("E6").Style

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Steve"
$ws.Range("C8").Value = "Carell"
$ws.Range("D8").Value = "testaccount7@mail.com"
$ws.Range("E8").Value = "Tester123@"
$ws.Range("F8").Value = "N"
$ws.Range("G8").Value = "N"
$ws.Range("A8").Style = $ws.Range("A6").Style
$ws.Range("E8").Style = $ws.Range("E6").Style

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Mark"
$ws.Range("C9").Value = "Doe"
$ws.Range("D9").Value = "testaccount8@mail.com"
$ws.Range("E9").Value = "Tester123@"
$ws.Range("F9").Value = "N"
$ws.Range("G9").Value = "N"
$ws.Range("A9").Style = $ws.Range("A6").Style
$ws.Range("E9").Style = $ws.Range("E6").Style

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Mark"
$ws.Range("C10").Value = "Doe"
$ws.Range("D10").Value = "testaccount9@mail.com"
$ws.Range("E10").Value = "Tester123@"
$ws.Range("F10").Value = "Y"
$ws.Range("G10").Value = "N"
$ws.Range("A10").Style = $ws.Range("A6").Style
$ws.Range("E10").Style = $ws.Range("E6").Style

# ---------------------------------------------------------------------------
# 3. Hyperlinks for the username/email column - rebuild all of them so the
#    tooltip is consistently present on every row (including row 2, which
#    previously lacked one) and the three new rows get theirs too.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:testaccount1@mail.com", "", "mailto:testaccount1@mail.com", "testaccount1@mail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:testaccount2@mail.com", "", "mailto:testaccount2@mail.com", "testaccount2@mail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:testaccount3@mail.com", "", "mailto:testaccount3@mail.com", "testaccount3@mail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:testaccount4@mail.com", "", "mailto:testaccount4@mail.com", "testaccount4@mail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:testaccount5@mail.com", "", "mailto:testaccount5@mail.com", "testaccount5@mail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:testaccount6@mail.com", "", "mailto:testaccount6@mail.com", "testaccount6@mail.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:testaccount7@mail.com", "", "mailto:testaccount7@mail.com", "testaccount7@mail.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:testaccount8@mail.com", "", "mailto:testaccount8@mail.com", "testaccount8@mail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:testaccount9@mail.com", "", "mailto:testaccount9@mail.com", "testaccount9@mail.com")

# ---------------------------------------------------------------------------
# 4. View / selection tweaks on the user sheet.
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()

Write-Output "user sheet updated"
